# Update cryptos list: prices (column D) and volume/1h change percentages (column E)
# Values are text strings in the original workbook (t="inlineStr"), including
# numeric-looking prices (e.g. "567.82") and dotted/thousands-style prices
# (e.g. "62.910.70"). Excel's COM Value setter auto-detects numeric-looking
# strings and converts them to real numbers, so for those cells we force the
# cell format to Text before assigning, then restore the default ("Normal")
# style afterwards so no stray style index gets baked into the cell (matching
# the original, which has no explicit style on these data cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = "D"; Value = "62.910.70"; ForceText = $false },
    @{ Row = 2; Col = "E"; Value = "  +1.59%  "; ForceText = $false },
    @{ Row = 3; Col = "D"; Value = "2.445.48"; ForceText = $false },
    @{ Row = 3; Col = "E"; Value = "  +1.87%  "; ForceText = $false },
    @{ Row = 4; Col = "E"; Value = "  +0.10%  "; ForceText = $false },
    @{ Row = 5; Col = "D"; Value = "567.82"; ForceText = $true },
    @{ Row = 5; Col = "E"; Value = "  +1.11%  "; ForceText = $false },
    @{ Row = 6; Col = "D"; Value = "146.39"; ForceText = $true },
    @{ Row = 6; Col = "E"; Value = "  +2.77%  "; ForceText = $false },
    @{ Row = 7; Col = "D"; Value = "0.999"; ForceText = $true },
    @{ Row = 7; Col = "E"; Value = "  -0.08%  "; ForceText = $false },
    @{ Row = 8; Col = "D"; Value = "0.535"; ForceText = $true },
    @{ Row = 8; Col = "E"; Value = "  +0.14%  "; ForceText = $false },
    @{ Row = 9; Col = "E"; Value = "  +2.93%  "; ForceText = $false },
    @{ Row = 10; Col = "E"; Value = "  +0.44%  "; ForceText = $false },
    @{ Row = 11; Col = "D"; Value = "5.31"; ForceText = $true },
    @{ Row = 11; Col = "E"; Value = "  +0.97%  "; ForceText = $false },
    @{ Row = 12; Col = "D"; Value = "0.356"; ForceText = $true },
    @{ Row = 12; Col = "E"; Value = "  +2.08%  "; ForceText = $false },
    @{ Row = 13; Col = "D"; Value = "27.07"; ForceText = $true },
    @{ Row = 13; Col = "E"; Value = "  +5.92%  "; ForceText = $false },
    @{ Row = 14; Col = "E"; Value = "  +6.49%  "; ForceText = $false },
    @{ Row = 15; Col = "D"; Value = "2.799.54"; ForceText = $false },
    @{ Row = 15; Col = "E"; Value = "  -1.18%  "; ForceText = $false },
    @{ Row = 16; Col = "D"; Value = "62.690.54"; ForceText = $false },
    @{ Row = 16; Col = "E"; Value = "  +1.44%  "; ForceText = $false },
    @{ Row = 17; Col = "D"; Value = "2.438.41"; ForceText = $false },
    @{ Row = 17; Col = "E"; Value = "  +1.57%  "; ForceText = $false },
    @{ Row = 18; Col = "D"; Value = "11.30"; ForceText = $true },
    @{ Row = 18; Col = "E"; Value = "  +0.85%  "; ForceText = $false },
    @{ Row = 19; Col = "D"; Value = "6.96"; ForceText = $true },
    @{ Row = 19; Col = "E"; Value = "  +2.39%  "; ForceText = $false },
    @{ Row = 20; Col = "D"; Value = "324.51"; ForceText = $true },
    @{ Row = 20; Col = "E"; Value = "  +1.22%  "; ForceText = $false },
    @{ Row = 21; Col = "E"; Value = "  +1.26%  "; ForceText = $false },
    @{ Row = 22; Col = "D"; Value = "0.999"; ForceText = $true },
    @{ Row = 22; Col = "E"; Value = "  -0.12%  "; ForceText = $false },
    @{ Row = 23; Col = "E"; Value = "  +7.34%  "; ForceText = $false },
    @{ Row = 24; Col = "D"; Value = "67.41"; ForceText = $true },
    @{ Row = 24; Col = "E"; Value = "  +2.22%  "; ForceText = $false },
    @{ Row = 25; Col = "D"; Value = "8.71"; ForceText = $true },
    @{ Row = 25; Col = "E"; Value = "  -0.75%  "; ForceText = $false },
    @{ Row = 26; Col = "D"; Value = "591.30"; ForceText = $true },
    @{ Row = 26; Col = "E"; Value = "  +5.28%  "; ForceText = $false },
    @{ Row = 27; Col = "E"; Value = "  +10.05%  "; ForceText = $false },
    @{ Row = 28; Col = "D"; Value = "2.565.26"; ForceText = $false },
    @{ Row = 28; Col = "E"; Value = "  +1.79%  "; ForceText = $false },
    @{ Row = 29; Col = "D"; Value = "8.48"; ForceText = $true },
    @{ Row = 29; Col = "E"; Value = "  +4.05%  "; ForceText = $false },
    @{ Row = 30; Col = "D"; Value = "0.998"; ForceText = $true },
    @{ Row = 30; Col = "E"; Value = "  -0.26%  "; ForceText = $false },
    @{ Row = 31; Col = "E"; Value = "  +5.43%  "; ForceText = $false },
    @{ Row = 32; Col = "E"; Value = "  +0.75%  "; ForceText = $false },
    @{ Row = 33; Col = "E"; Value = "  +0.79%  "; ForceText = $false },
    @{ Row = 34; Col = "E"; Value = "  +3.19%  "; ForceText = $false },
    @{ Row = 35; Col = "D"; Value = "4.89"; ForceText = $true },
    @{ Row = 35; Col = "E"; Value = "  +4.33%  "; ForceText = $false },
    @{ Row = 36; Col = "D"; Value = "0.998"; ForceText = $true },
    @{ Row = 36; Col = "E"; Value = "  -0.17%  "; ForceText = $false },
    @{ Row = 37; Col = "E"; Value = "  +1.57%  "; ForceText = $false },
    @{ Row = 38; Col = "E"; Value = "  +0.77%  "; ForceText = $false },
    @{ Row = 39; Col = "D"; Value = "18.84"; ForceText = $true },
    @{ Row = 39; Col = "E"; Value = "  +1.52%  "; ForceText = $false },
    @{ Row = 40; Col = "D"; Value = "148.74"; ForceText = $true },
    @{ Row = 40; Col = "E"; Value = "  -2.31%  "; ForceText = $false },
    @{ Row = 41; Col = "E"; Value = "  +2.74%  "; ForceText = $false },
    @{ Row = 43; Col = "D"; Value = "2.47"; ForceText = $true },
    @{ Row = 43; Col = "E"; Value = "  +10.23%  "; ForceText = $false },
    @{ Row = 44; Col = "D"; Value = "149.49"; ForceText = $true },
    @{ Row = 44; Col = "E"; Value = "  +1.46%  "; ForceText = $false },
    @{ Row = 45; Col = "D"; Value = "3.70"; ForceText = $true },
    @{ Row = 45; Col = "E"; Value = "  +2.75%  "; ForceText = $false },
    @{ Row = 46; Col = "D"; Value = "0.0538"; ForceText = $true },
    @{ Row = 46; Col = "E"; Value = "  +1.80%  "; ForceText = $false },
    @{ Row = 47; Col = "D"; Value = "20.69"; ForceText = $true },
    @{ Row = 47; Col = "E"; Value = "  +4.79%  "; ForceText = $false },
    @{ Row = 48; Col = "E"; Value = "  +3.08%  "; ForceText = $false },
    @{ Row = 49; Col = "E"; Value = "  +3.67%  "; ForceText = $false },
    @{ Row = 50; Col = "D"; Value = "0.0925"; ForceText = $true },
    @{ Row = 50; Col = "E"; Value = "  +1.03%  "; ForceText = $false },
    @{ Row = 51; Col = "E"; Value = "  +4.35%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
